# Rewrite the missing-data table (A2:F33) to match the post-edit state.
# The edit re-imputes / blanks several scattered cells in rows 2-25, and
# shifts rows 26-35 up by one (dropping the old row 26 "RM 232" entry and
# the trailing rows 34-35), with several value corrections along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('RM 2', -19.7, 14.9, -13.5, -7.2, 18.03),
    @('RM 8', -19.7, 11.2, -14.2, -5.7, 17.64),
    @('RM 9', -18.7, 11, -15.4, $null, 17.97),
    @('RM 14', -19.5, 12.3, $null, -5, 17.66),
    @('RM 21', -19.8, 15.1, -14.2, -5.7, 16.43),
    @('RM 32', -19.5, 15, -13.8, -7.1, 17.24),
    @('RM 38', -19.9, 15.5, -13.9, -6.6, 17.05),
    @('RM 42', -20.6, 10.5, -14.5, -6.8, 17.26),
    @('RM 52 a', $null, 11.5, -14.7, -6.1, 16.43),
    @('RM 58', -20.8, 11.4, -15.5, -7.9, 17.65),
    @('RM 81', $null, 12.5, -14.1, -5.3, $null),
    @('RM 88', -19.9, 12.5, -13.9, $null, 17.1),
    @('RM 90', -19.6, 14.4, -13.1, $null, 17.76),
    @('RM 95', -19.1, 12.5, -15.2, -8.4, 16.2),
    @('RM 103', -19.5, 13.5, -15.3, -6.9, 17.34),
    @('RM 116', -19.4, 11.2, -14.7, -7.3, 17.78),
    @('RM 120', $null, 11.5, -15.2, -8.5, 18.35),
    @('RM 125', -20.6, 13.2, -15.5, -6.5, 17.81),
    @('RM 134', -19, 12.5, -14, -7.2, $null),
    @('RM 135', -18.9, 12.7, -14.3, -8.7, 16.58),
    @('RM 138', -19.3, 13.5, -15.4, -6.1, 16.81),
    @('RM 140', -19.5, 12.2, $null, -7, 16.48),
    @('RM 142a', -17.7, 12.7, -13.9, -8.1, 16.78),
    @('RM 145', -19.5, 10.7, -15.5, -7.1, $null),
    @('SC 5', -20.2, 10.8, -13.8, -5, 17.38),
    @('SC 101', -20.4, 10, -14.6, -10, 17),
    @('SC 105', -19.6, 11.1, -13.7, -5.9, 17.44),
    @('SC 119', -19.5, 11.2, -13, -6.8, 18.06),
    @('SC 120', -19.7, 11.4, -13.6, -5.7, 16.89),
    @('SC 132', -18.8, 15.3, -13.7, -8.1, 17.18),
    @('SC 193', -19.9, 10.5, -14.7, -6.4, 17.39),
    @('SC 232', -19.5, 10.4, -14.1, -10.7, 17.53)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    for ($c = 1; $c -lt $row.Count; $c++) {
        $colLetter = [char](65 + $c)
        $cell = $ws.Range("$colLetter$r")
        if ($null -eq $row[$c]) {
            $cell.Value = $null
        } else {
            $cell.Value = $row[$c]
        }
    }
}

# Rows 34 and 35 (old "SC 193" / "SC 232") no longer exist — clear them so
# the sheet's used range shrinks back down to A1:F33.
$ws.Range("A34:F35").Clear()
